$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ridership")

# Update Riders (column C) and Average (column D) values for the Madigan bike hours update
$ws.Range("C2").Value = 257
$ws.Range("D2").Value = 263

$ws.Range("C3").Value = 209
$ws.Range("D3").Value = 234.5

$ws.Range("C4").Value = 232
$ws.Range("D4").Value = 251

$ws.Range("C5").Value = 221
$ws.Range("D5").Value = 243.5

$ws.Range("C6").Value = 116
$ws.Range("D6").Value = 113

$ws.Range("C7").Value = 66
$ws.Range("D7").Value = 74

$wb.Save()
